$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.007.47'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '2.227.28'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '294.39'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D6').Value = '86.84'
$ws.Range('E6').Value = '  -1.29%  '
$ws.Range('D7').Value = '0.512'
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.467'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').Value = '30.48'
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '50.61'
$ws.Range('E11').Value = '  +6.27%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = '0.0781'
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('D14').Value = '6.42'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').Value = '2.566.52'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '13.83'
$ws.Range('E16').Value = '  -1.42%  '
$ws.Range('D17').Value = '2.222.54'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('D18').Value = '0.735'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').Value = '39.919.28'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').Value = '0.0₃0885'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').Value = '11.20'
$ws.Range('E21').Value = '  -5.03%  '
$ws.Range('D22').Value = '5.77'
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('D23').Value = '65.73'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '235.78'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').Value = '2.47'
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('E28').Value = '  +6.15%  '
$ws.Range('D29').Value = '23.18'
$ws.Range('E29').Value = '  +2.42%  '
$ws.Range('D30').Value = '9.28'
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('D31').Value = '157.86'
$ws.Range('E31').Value = '  +3.61%  '
$ws.Range('D32').Value = '31.50'
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = '4.95'
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').Value = '3.01'
$ws.Range('E35').Value = '  +5.60%  '
$ws.Range('D36').Value = '0.0712'
$ws.Range('E36').Value = '  -1.22%  '
$ws.Range('D37').Value = '2.30'
$ws.Range('E37').Value = '  -3.43%  '
$ws.Range('D38').Value = '0.113'
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('D39').Value = '0.0987'
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('D40').Value = '1.74'
$ws.Range('E40').Value = '  +1.69%  '
$ws.Range('D41').Value = '15.38'
$ws.Range('E41').Value = '  -4.65%  '
$ws.Range('D42').Value = '2.087.58'
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('E43').Value = '  -3.74%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '18.09'
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('D45').Value = '10.07'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0269'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('D47').Value = '1.96'
$ws.Range('E47').Value = '  -9.70%  '
$ws.Range('D48').Value = '2.69'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('D49').Value = '2.439.12'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.46'
$ws.Range('E50').Value = '  +0.96%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '1.11'
$ws.Range('E51').Value = '  +3.23%  '
